# Update relay settings (F2 relay - Banshee) with back-of-a-napkin calculated
# values for CT Primary (G), 51P TOC Trip Pickup (K), 27P Trip Pickup (N),
# and 59P Trip Pickup (O) across relay rows 2-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1200
$ws.Range("K2").Value = 3.4135813600798519
$ws.Range("N2").Value = 0.8
$ws.Range("O2").Value = 1.2
$ws.Range("G3").Value = 1200
$ws.Range("K3").Value = 1.7039851694832704
$ws.Range("N3").Value = 0.8
$ws.Range("O3").Value = 1.2
$ws.Range("G4").Value = 1200
$ws.Range("K4").Value = 0.80187537387448016
$ws.Range("N4").Value = 0.8
$ws.Range("O4").Value = 1.2
$ws.Range("G5").Value = 1200
$ws.Range("K5").Value = 0.99752525475971998
$ws.Range("N5").Value = 0.8
$ws.Range("O5").Value = 1.2
$ws.Range("G6").Value = 1200
$ws.Range("K6").Value = 0.50117210867155004
$ws.Range("N6").Value = 0.8
$ws.Range("O6").Value = 1.2
$ws.Range("G7").Value = 1200
$ws.Range("K7").Value = 0.10023442173431002
$ws.Range("N7").Value = 0.8
$ws.Range("O7").Value = 1.2
$ws.Range("G8").Value = 1200
$ws.Range("K8").Value = 0.75175816300732523
$ws.Range("N8").Value = 0.8
$ws.Range("O8").Value = 1.2
$ws.Range("G9").Value = 1200
$ws.Range("K9").Value = 0.66501683650647991
$ws.Range("N9").Value = 0.8
$ws.Range("O9").Value = 1.2
$ws.Range("G10").Value = 1200
$ws.Range("K10").Value = 0.50117210867155004
$ws.Range("N10").Value = 0.8
$ws.Range("O10").Value = 1.2
$ws.Range("G11").Value = 1200
$ws.Range("K11").Value = 0.20046884346862004
$ws.Range("N11").Value = 0.8
$ws.Range("O11").Value = 1.2
$ws.Range("G12").Value = 1200
$ws.Range("K12").Value = 1.0023442173431001
$ws.Range("N12").Value = 0.8
$ws.Range("O12").Value = 1.2
$ws.Range("G13").Value = 3500
$ws.Range("K13").Value = 0.79042001139055906
$ws.Range("N13").Value = 0.8
$ws.Range("O13").Value = 1.2
$ws.Range("G14").Value = 3500
$ws.Range("K14").Value = 1.4128571428571426
$ws.Range("N14").Value = 0.8
$ws.Range("O14").Value = 1.2
$ws.Range("G15").Value = 1200
$ws.Range("K15").Value = 0.40093768693724008
$ws.Range("N15").Value = 0.8
$ws.Range("O15").Value = 1.2
$ws.Range("G16").Value = 1200
$ws.Range("K16").Value = 0.40093768693724008
$ws.Range("N16").Value = 0.8
$ws.Range("O16").Value = 1.2
$ws.Range("G17").Value = 3500
$ws.Range("K17").Value = 1.4128571428571426
$ws.Range("N17").Value = 0.8
$ws.Range("O17").Value = 1.2
$ws.Range("G18").Value = 3500
$ws.Range("K18").Value = 1.4128571428571426
$ws.Range("N18").Value = 0.8
$ws.Range("O18").Value = 1.2
$ws.Range("G19").Value = 3500
$ws.Range("K19").Value = 1.3832350199334786
$ws.Range("N19").Value = 0.8
$ws.Range("O19").Value = 1.2
$ws.Range("G20").Value = 3500
$ws.Range("K20").Value = 1.976050028476398
$ws.Range("N20").Value = 0.8
$ws.Range("O20").Value = 1.2

# The 51P TOC Trip Pickup column now carries calculated (non-round) values,
# so format it to show two decimal places.
$ws.Range("K2:K20").NumberFormat = "0.00"

# Leave the selection where the author left off editing.
$ws.Range("O23").Select() | Out-Null
